# Add 2021 (remainder) and 2022 ICS values, plus start of 2023, to the
# MeadLevelEndOfMonth table.
#
# xlPasteFormats = -4122 (used to copy cell formatting/styles without values)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Prime the formatting for the two brand-new rows (89 = 2022, 90 =
#    2023) and row 88's currently-blank trailing cells BEFORE we start
#    touching row 88's own styles, so we can still copy its existing
#    "empty" styles (s=5 / s=6) as a template.
# ---------------------------------------------------------------------

# Row 87 (2020) has the canonical fully-filled-in style pattern
# (col A bold/bordered = style 3, cols B:M = style 4). Stamp that
# pattern onto the two new rows 89 and 90.
$ws.Range("A87:M87").Copy()
$ws.Range("A89:M89").PasteSpecial(-4122)
$ws.Range("A90:M90").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 90 only has data through column H; I90 stays style 4 (blank),
# J90:L90 need the "blank" style 5, and M90 the "blank" style 6 - copy
# those straight from row 88's (still untouched) trailing blank cells.
$ws.Range("G88:I88").Copy()
$ws.Range("J90:L90").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("M88").Copy()
$ws.Range("M90").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Row 88 (2021): fill in the rest of the year (G:M) using the
#    regular data style (style 4, same as F88).
# ---------------------------------------------------------------------
$ws.Range("F88").Copy()
$ws.Range("G88:M88").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(88, 7).Value  = 1068.77
$ws.Cells.Item(88, 8).Value  = 1067.6500000000001
$ws.Cells.Item(88, 9).Value  = 1067.96
$ws.Cells.Item(88, 10).Value = 1067.68
$ws.Cells.Item(88, 11).Value = 1066.77
$ws.Cells.Item(88, 12).Value = 1064.97
$ws.Cells.Item(88, 13).Value = 1066.3900000000001

# ---------------------------------------------------------------------
# 3) Row 89 (2022): full year of data.
# ---------------------------------------------------------------------
$ws.Cells.Item(89, 1).Value  = 2022
$ws.Cells.Item(89, 2).Value  = 1067.0899999999999
$ws.Cells.Item(89, 3).Value  = 1066.78
$ws.Cells.Item(89, 4).Value  = 1061.49
$ws.Cells.Item(89, 5).Value  = 1054.69
$ws.Cells.Item(89, 6).Value  = 1047.69
$ws.Cells.Item(89, 7).Value  = 1043.02
$ws.Cells.Item(89, 8).Value  = 1040.92
$ws.Cells.Item(89, 9).Value  = 1044.28
$ws.Cells.Item(89, 10).Value = 1045.03
$ws.Cells.Item(89, 11).Value = 1046.28
$ws.Cells.Item(89, 12).Value = 1043.02
$ws.Cells.Item(89, 13).Value = 1044.82

# ---------------------------------------------------------------------
# 4) Row 90 (2023): partial year, through August (column H) only.
# ---------------------------------------------------------------------
$ws.Cells.Item(90, 1).Value = 2023
$ws.Cells.Item(90, 2).Value = 1046.97
$ws.Cells.Item(90, 3).Value = 1047.02
$ws.Cells.Item(90, 4).Value = 1046.03
$ws.Cells.Item(90, 5).Value = 1049.69
$ws.Cells.Item(90, 6).Value = 1054.28
$ws.Cells.Item(90, 7).Value = 1056.3900000000001
$ws.Cells.Item(90, 8).Value = 1061.02

# ---------------------------------------------------------------------
# 5) Row 88, column A: the "2021" year label becomes a text value
#    (shared string) rather than a number, matching how the other
#    year labels above it (e.g. "2019" in A86) are already stored.
#    Forcing text requires a Text number format while the value is
#    entered; re-apply the original (style 3) formatting afterwards
#    so the visible formatting is unchanged.
# ---------------------------------------------------------------------
$ws.Range("A88").NumberFormat = "@"
$ws.Range("A88").Value = "2021"
$ws.Range("A87").Copy()
$ws.Range("A88").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 6) Leave the selection where the editor last left off.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("M103").Select()
